$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(41, 8).Value = 150.5
$ws.Cells.Item(41, 9).Value = 150.5
$ws.Cells.Item(41, 11).Value = 150.5
$ws.Cells.Item(41, 13).Value = 289.5

$ws.Cells.Item(113, 8).Value = 244275.44
$ws.Cells.Item(113, 9).Value = 373065.66
$ws.Cells.Item(113, 10).Value = 2793.75
$ws.Cells.Item(113, 11).Value = 373065.66
$ws.Cells.Item(113, 12).Value = 2793.75
$ws.Cells.Item(113, 13).Value = -369811.66
$ws.Cells.Item(113, 14).Value = -9301.75

$ws.Cells.Item(135, 8).Value = 2603.5557
$ws.Cells.Item(135, 9).Value = 2761.923
$ws.Cells.Item(135, 10).Value = 2191.8
$ws.Cells.Item(135, 11).Value = 24857.307
$ws.Cells.Item(135, 12).Value = 19726.2
$ws.Cells.Item(135, 13).Value = -22322.307
$ws.Cells.Item(135, 14).Value = -24796.2

$ws.Cells.Item(137, 8).Value = 974.8387
$ws.Cells.Item(137, 9).Value = 802.5961
$ws.Cells.Item(137, 10).Value = 1870.5
$ws.Cells.Item(137, 11).Value = 2407.7883
$ws.Cells.Item(137, 12).Value = 5611.5
$ws.Cells.Item(137, 13).Value = 142.2116999999998
$ws.Cells.Item(137, 14).Value = -10711.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 1032.9166
$ws.Cells.Item(74, 9).Value = 1010.0909
$ws.Cells.Item(74, 10).Value = 1052.2307
$ws.Cells.Item(74, 11).Value = 1010.0909
$ws.Cells.Item(74, 12).Value = 1052.2307
$ws.Cells.Item(74, 13).Value = -136.0909
$ws.Cells.Item(74, 14).Value = -2800.2307

$ws.Cells.Item(77, 8).Value = 1032.9166
$ws.Cells.Item(77, 9).Value = 1010.0909
$ws.Cells.Item(77, 10).Value = 1052.2307
$ws.Cells.Item(77, 11).Value = 5050.4545
$ws.Cells.Item(77, 12).Value = 5261.1535
$ws.Cells.Item(77, 13).Value = -682.4544999999998
$ws.Cells.Item(77, 14).Value = -13997.1535

$ws.Cells.Item(132, 8).Value = 1046.459
$ws.Cells.Item(132, 9).Value = 835.4400000000001
$ws.Cells.Item(132, 10).Value = 2005.6364
$ws.Cells.Item(132, 11).Value = 2506.32
$ws.Cells.Item(132, 12).Value = 6016.9092
$ws.Cells.Item(132, 13).Value = 23.67999999999984
$ws.Cells.Item(132, 14).Value = -11076.9092

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 910.4400000000001
$ws.Cells.Item(94, 9).Value = 884.3333
$ws.Cells.Item(94, 10).Value = 1047.5
$ws.Cells.Item(94, 11).Value = 884.3333
$ws.Cells.Item(94, 12).Value = 1047.5
$ws.Cells.Item(94, 13).Value = -433.3333
$ws.Cells.Item(94, 14).Value = -1949.5

$ws.Cells.Item(134, 8).Value = 2085.6924
$ws.Cells.Item(134, 9).Value = 1884
$ws.Cells.Item(134, 10).Value = 2233.6
$ws.Cells.Item(134, 11).Value = 5652
$ws.Cells.Item(134, 12).Value = 6700.799999999999
$ws.Cells.Item(134, 13).Value = -3117
$ws.Cells.Item(134, 14).Value = -11770.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 780.39624
$ws.Cells.Item(58, 9).Value = 696.1842
$ws.Cells.Item(58, 10).Value = 993.73334
$ws.Cells.Item(58, 11).Value = 696.1842
$ws.Cells.Item(58, 12).Value = 993.73334
$ws.Cells.Item(58, 13).Value = -493.1842
$ws.Cells.Item(58, 14).Value = -1399.73334

$ws.Cells.Item(132, 8).Value = 1367.7317
$ws.Cells.Item(132, 9).Value = 1049.3823
$ws.Cells.Item(132, 10).Value = 2914
$ws.Cells.Item(132, 11).Value = 3148.1469
$ws.Cells.Item(132, 12).Value = 8742
$ws.Cells.Item(132, 13).Value = -618.1468999999997
$ws.Cells.Item(132, 14).Value = -13802

$ws.Cells.Item(134, 8).Value = 1291.1111
$ws.Cells.Item(134, 9).Value = 1093.3334
$ws.Cells.Item(134, 11).Value = 3280.0002
$ws.Cells.Item(134, 13).Value = -745.0001999999999

$ws.Cells.Item(136, 8).Value = 780.39624
$ws.Cells.Item(136, 9).Value = 696.1842
$ws.Cells.Item(136, 10).Value = 993.73334
$ws.Cells.Item(136, 11).Value = 2088.5526
$ws.Cells.Item(136, 12).Value = 2981.20002
$ws.Cells.Item(136, 13).Value = 461.4474
$ws.Cells.Item(136, 14).Value = -8081.20002

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 365.78125
$ws.Cells.Item(5, 9).Value = 261.73077
$ws.Cells.Item(5, 10).Value = 816.6667
$ws.Cells.Item(5, 11).Value = 785.19231
$ws.Cells.Item(5, 12).Value = 2450.0001
$ws.Cells.Item(5, 13).Value = -673.19231
$ws.Cells.Item(5, 14).Value = -2674.0001

$ws.Cells.Item(68, 8).Value = 10590.6
$ws.Cells.Item(68, 9).Value = 33734.668
$ws.Cells.Item(68, 10).Value = 671.7143
$ws.Cells.Item(68, 11).Value = 101204.004
$ws.Cells.Item(68, 12).Value = 2015.1429
$ws.Cells.Item(68, 13).Value = -100393.004
$ws.Cells.Item(68, 14).Value = -3637.1429

$ws.Cells.Item(70, 8).Value = 1651.8182
$ws.Cells.Item(70, 9).Value = 1146.25
$ws.Cells.Item(70, 10).Value = 3000
$ws.Cells.Item(70, 11).Value = 3438.75
$ws.Cells.Item(70, 12).Value = 9000
$ws.Cells.Item(70, 13).Value = -3123.75
$ws.Cells.Item(70, 14).Value = -9630

$ws.Cells.Item(71, 8).Value = 10590.6
$ws.Cells.Item(71, 9).Value = 33734.668
$ws.Cells.Item(71, 10).Value = 671.7143
$ws.Cells.Item(71, 11).Value = 303612.012
$ws.Cells.Item(71, 12).Value = 6045.428699999999
$ws.Cells.Item(71, 13).Value = -299556.012
$ws.Cells.Item(71, 14).Value = -14157.4287

$ws.Cells.Item(73, 8).Value = 1651.8182
$ws.Cells.Item(73, 9).Value = 1146.25
$ws.Cells.Item(73, 10).Value = 3000
$ws.Cells.Item(73, 11).Value = 3438.75
$ws.Cells.Item(73, 12).Value = 9000
$ws.Cells.Item(73, 13).Value = -2346.75
$ws.Cells.Item(73, 14).Value = -11184

$ws.Cells.Item(92, 8).Value = 889.7368
$ws.Cells.Item(92, 9).Value = 464.57144
$ws.Cells.Item(92, 10).Value = 1137.75
$ws.Cells.Item(92, 11).Value = 1393.71432
$ws.Cells.Item(92, 12).Value = 3413.25
$ws.Cells.Item(92, 13).Value = -145.71432
$ws.Cells.Item(92, 14).Value = -5909.25

$ws.Cells.Item(131, 8).Value = 863.60657
$ws.Cells.Item(131, 10).Value = 967.4400000000001
$ws.Cells.Item(131, 12).Value = 2902.32
$ws.Cells.Item(131, 14).Value = -12982.32

$ws.Cells.Item(135, 8).Value = 365.78125
$ws.Cells.Item(135, 9).Value = 261.73077
$ws.Cells.Item(135, 10).Value = 816.6667
$ws.Cells.Item(135, 11).Value = 2355.57693
$ws.Cells.Item(135, 12).Value = 7350.0003
$ws.Cells.Item(135, 13).Value = 179.4230699999998
$ws.Cells.Item(135, 14).Value = -12420.0003

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(68, 8).Value = 30000
$ws.Cells.Item(68, 10).Value = 30000
$ws.Cells.Item(68, 12).Value = 30000
$ws.Cells.Item(68, 14).Value = -31622

$ws.Cells.Item(71, 8).Value = 30000
$ws.Cells.Item(71, 10).Value = 30000
$ws.Cells.Item(71, 12).Value = 90000
$ws.Cells.Item(71, 14).Value = -98112

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 4483.25
$ws.Cells.Item(7, 9).Value = 4285.5713
$ws.Cells.Item(7, 10).Value = 4760
$ws.Cells.Item(7, 11).Value = 4285.5713
$ws.Cells.Item(7, 12).Value = 4760
$ws.Cells.Item(7, 13).Value = -4173.5713
$ws.Cells.Item(7, 14).Value = -4984

$ws.Cells.Item(16, 8).Value = 3901.7368
$ws.Cells.Item(16, 9).Value = 3268.875
$ws.Cells.Item(16, 10).Value = 4362
$ws.Cells.Item(16, 11).Value = 3268.875
$ws.Cells.Item(16, 12).Value = 4362
$ws.Cells.Item(16, 13).Value = -3098.875
$ws.Cells.Item(16, 14).Value = -4702

$ws.Cells.Item(46, 8).Value = 856.8333
$ws.Cells.Item(46, 9).Value = 627.8
$ws.Cells.Item(46, 10).Value = 2002
$ws.Cells.Item(46, 11).Value = 627.8
$ws.Cells.Item(46, 12).Value = 2002
$ws.Cells.Item(46, 13).Value = -439.8
$ws.Cells.Item(46, 14).Value = -2378

$ws.Cells.Item(55, 8).Value = 187.35484
$ws.Cells.Item(55, 9).Value = 164.57143
$ws.Cells.Item(55, 10).Value = 400
$ws.Cells.Item(55, 11).Value = 164.57143
$ws.Cells.Item(55, 12).Value = 400
$ws.Cells.Item(55, 13).Value = 8.428570000000008
$ws.Cells.Item(55, 14).Value = -746

$ws.Cells.Item(126, 8).Value = 4483.25
$ws.Cells.Item(126, 9).Value = 4285.5713
$ws.Cells.Item(126, 10).Value = 4760
$ws.Cells.Item(126, 11).Value = 12856.7139
$ws.Cells.Item(126, 12).Value = 14280
$ws.Cells.Item(126, 13).Value = -10386.7139
$ws.Cells.Item(126, 14).Value = -19220

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 1079.4166
$ws.Cells.Item(136, 9).Value = 939.5925999999999
$ws.Cells.Item(136, 10).Value = 1259.1904
$ws.Cells.Item(136, 11).Value = 2818.7778
$ws.Cells.Item(136, 12).Value = 3777.5712
$ws.Cells.Item(136, 13).Value = -268.7777999999998
$ws.Cells.Item(136, 14).Value = -8877.5712
